## "Haflway through a few more bits and pieces getting the networks together."
##
## Sheet1 tracks, per IEEE test circuit, which modelling pieces are done
## (Fixed model / Tap Model / LTC Model / HC Calcs columns L:O, marked with
## "X"/"O") plus a scratch "TODO ORDER" / "DIFF" note (columns P:Q) for the
## circuits that are not finished yet.
##
## This pass finishes off the Ckt24 row (row 16): the Tap/LTC/HC-calc
## columns are now implemented (M16/N16/O16 go from blank to X/O/X, matching
## the pattern already used on the rows above/below it), so the leftover
## "MED" / order-note scratch cell for that row is cleared out. The other
## two outstanding rows (13 and 14) just get their "DIFF" counters bumped
## down by one as part of the same clean-up pass.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ckt5 / Ckt7 "DIFF" counters (col Q) — one step closer to done.
$ws.Range("Q13").Value = 1
$ws.Range("Q14").Value = 2

# Ckt24 (row 16): Tap Model / LTC Model / HC Calcs now implemented.
$ws.Range("M16").Value = "X"
$ws.Range("N16").Value = "O"
$ws.Range("O16").Value = "X"

# ... so the "TODO ORDER"/"DIFF" scratch note for that row is no longer needed.
$ws.Range("P16:Q16").ClearContents()

# Leave the cursor where the author was last working.
$ws.Range("K9").Select()
